# Generate Report for Handoff
# Swaps the handoff/handback report rows for 39c3eb9d... and 6220edcd...
# so that 6220edcd moves to the top (row 2) and 39c3eb9d moves to row 3,
# updates statuses/timestamps, and records a "stale handback" error for
# the 39c3eb9d file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-27 20:47:41"

foreach ($hl in $ov.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = 'e2e\6220edcd-2394-46a0-ad90-271fd6f9c4b5.md'
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = 'e2e\39c3eb9d-8082-41dd-bd25-8260aac71e11.md'
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "6220edcd-2394-46a0-ad90-271fd6f9c4b5.md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("G2").Value = "6220edcd-2394-46a0-ad90-271fd6f9c4b5.7b3b161129b4923eb4dd032b6a4d99b4687c4630.zh-cn.xlf"
$zh.Range("I2").Value = "6220edcd-2394-46a0-ad90-271fd6f9c4b5.md"
$zh.Range("J2").Value = "6220edcd-2394-46a0-ad90-271fd6f9c4b5.7b3b161129b4923eb4dd032b6a4d99b4687c4630.zh-cn.xlf"

$zh.Range("A3").Value = "39c3eb9d-8082-41dd-bd25-8260aac71e11.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "39c3eb9d-8082-41dd-bd25-8260aac71e11.2648b9e340541eca31485e5434bd3a04b5f88a40.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-27 20:47:36"
$zh.Range("I3").Value = "39c3eb9d-8082-41dd-bd25-8260aac71e11.md"
$zh.Range("J3").Value = "39c3eb9d-8082-41dd-bd25-8260aac71e11.2648b9e340541eca31485e5434bd3a04b5f88a40.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d4732e5d32d1d0c0747633d0e9ed0c821bd7f76/e2e/39c3eb9d-8082-41dd-bd25-8260aac71e11.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6048583f19c80418f32c38df4943bbc2cc7b71d8/e2e/39c3eb9d-8082-41dd-bd25-8260aac71e11.md."

foreach ($hl in $zh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = '6220edcd-2394-46a0-ad90-271fd6f9c4b5.md'
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = '6220edcd-2394-46a0-ad90-271fd6f9c4b5.md'
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = '39c3eb9d-8082-41dd-bd25-8260aac71e11.md'
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = '39c3eb9d-8082-41dd-bd25-8260aac71e11.md'
    }
}

$zh.Range("P1").ColumnWidth = 39.1875

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "6220edcd-2394-46a0-ad90-271fd6f9c4b5.md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("G2").Value = "6220edcd-2394-46a0-ad90-271fd6f9c4b5.7b3b161129b4923eb4dd032b6a4d99b4687c4630.de-de.xlf"
$de.Range("H2").Value = "2016-08-27 20:47:41"
$de.Range("I2").Value = "6220edcd-2394-46a0-ad90-271fd6f9c4b5.md"
$de.Range("J2").Value = "6220edcd-2394-46a0-ad90-271fd6f9c4b5.7b3b161129b4923eb4dd032b6a4d99b4687c4630.de-de.xlf"

$de.Range("A3").Value = "39c3eb9d-8082-41dd-bd25-8260aac71e11.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "39c3eb9d-8082-41dd-bd25-8260aac71e11.2648b9e340541eca31485e5434bd3a04b5f88a40.de-de.xlf"
$de.Range("H3").Value = "2016-08-27 20:47:41"
$de.Range("I3").Value = "39c3eb9d-8082-41dd-bd25-8260aac71e11.md"
$de.Range("J3").Value = "39c3eb9d-8082-41dd-bd25-8260aac71e11.2648b9e340541eca31485e5434bd3a04b5f88a40.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4d4732e5d32d1d0c0747633d0e9ed0c821bd7f76/e2e/39c3eb9d-8082-41dd-bd25-8260aac71e11.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6048583f19c80418f32c38df4943bbc2cc7b71d8/e2e/39c3eb9d-8082-41dd-bd25-8260aac71e11.md."

foreach ($hl in $de.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = '6220edcd-2394-46a0-ad90-271fd6f9c4b5.md'
    } elseif ($addr -eq '$I$2') {
        $hl.TextToDisplay = '6220edcd-2394-46a0-ad90-271fd6f9c4b5.md'
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = '39c3eb9d-8082-41dd-bd25-8260aac71e11.md'
    } elseif ($addr -eq '$I$3') {
        $hl.TextToDisplay = '39c3eb9d-8082-41dd-bd25-8260aac71e11.md'
    }
}

$de.Range("P1").ColumnWidth = 39.1875
